$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches source inlineStr string cells) for all edited
# cells so numeric-looking strings ("256.13") and percentages ("-0.09%") are
# not auto-converted to numbers by the input parser.
$cells = @{
    'D2' = '256.13'
    'E2' = '-0.09%'
    'D3' = '26.60'
    'E3' = '-1.63%'
    'D4' = '4.638'
    'E4' = '0.41%'
    'D5' = '0.05922'
    'E5' = '0.47%'
    'D6' = '6.602'
    'E6' = '-0.69%'
    'D7' = '0.8541'
    'E7' = '-1.60%'
    'D8' = '0.9103'
    'E8' = '-3.90%'
    'D9' = '0.1379'
    'E9' = '-1.67%'
    'D10' = '0.04327'
    'D11' = '0.06998'
    'E11' = '-1.26%'
    'D12' = '0.03019'
    'E12' = '-5.91%'
    'D13' = '0.09102'
    'D14' = '0.001531'
    'E14' = '-1.14%'
    'B15' = 'One'
    'C15' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'D15' = '0.0006060'
    'E15' = '0.61%'
    'B16' = 'TigerCash'
    'C16' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D16' = '0.006030'
    'E16' = '-0.21%'
    'B17' = 'LEO'
    'C17' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D17' = '3.475'
    'E17' = '-1.10%'
    'B18' = 'GateToken'
    'C18' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D18' = '3.129'
    'E18' = '-1.96%'
    'B19' = 'BTSEToken'
    'C19' = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
    'D19' = '2.150'
    'E19' = '-2.36%'
    'D20' = '0.3080'
    'E20' = '0.17%'
    'D21' = '0.1285'
    'E21' = '0.17%'
    'D22' = '3.889'
    'E22' = '0.95%'
    'D23' = '0.04211'
    'E23' = '-0.50%'
    'D24' = '0.001217'
    'E24' = '-0.26%'
    'D25' = '0.004626'
    'E25' = '8.24%'
    'D26' = '0.0001201'
    'E26' = '0.09%'
    'D27' = '0.0001712'
    'E27' = '13.49%'
    'D40' = '0.03799'
    'E40' = '-0.45%'
    'B41' = 'BKEXToken'
    'C41' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'D41' = '0.1101'
    'E41' = '0.10%'
    'B42' = 'KickToken'
    'C42' = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
    'D42' = '0.003686'
    'E42' = '-40.44%'
    'D43' = '0.002411'
    'E43' = '5.92%'
    'E44' = '23.83%'
    'D45' = '0.00005120'
    'E45' = '-6.96%'
    'D46' = '0.00000000749'
    'E46' = '-0.24%'
    'D47' = '0.04991'
    'E47' = '-17.14%'
    'E48' = '10,469.86%'
    'D49' = '0.00002096'
    'E49' = '-0.24%'
    'D50' = '0.0001996'
    'E50' = '-0.24%'
}

foreach ($ref in $cells.Keys) {
    $range = $ws.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $cells[$ref]
}
